$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "1.00", "603.27") are preserved verbatim as text, matching
# the source data which stores prices as inline strings.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: D2, E2
$ws.Range("D2").Value = '72.616.23'
$ws.Range("E2").Value = '  -0.60%  '

# Row 3: D3, E3
$ws.Range("D3").Value = '3.932.65'
$ws.Range("E3").Value = '  -2.76%  '

# Row 4: D4, E4
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '

# Row 5: D5, E5
$ws.Range("D5").Value = '603.27'
$ws.Range("E5").Value = '  +1.57%  '

# Row 6: D6, E6
$ws.Range("D6").Value = '170.93'
$ws.Range("E6").Value = '  +11.14%  '

# Row 7: D7, E7
$ws.Range("D7").Value = '0.683'
$ws.Range("E7").Value = '  -0.88%  '

# Row 8: D8, E8
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.05%  '

# Row 9: D9, E9
$ws.Range("D9").Value = '0.787'
$ws.Range("E9").Value = '  +3.19%  '

# Row 10: D10, E10
$ws.Range("D10").Value = '0.183'
$ws.Range("E10").Value = '  +7.00%  '

# Row 11: D11, E11
$ws.Range("D11").Value = '55.70'
$ws.Range("E11").Value = '  +2.98%  '

# Row 12: D12, E12
$ws.Range("D12").Value = '0.0000327'
$ws.Range("E12").Value = '  +1.23%  '

# Row 13: D13, E13
$ws.Range("D13").Value = '11.60'
$ws.Range("E13").Value = '  +5.26%  '

# Row 14: D14, E14
$ws.Range("D14").Value = '4.577.41'
$ws.Range("E14").Value = '  -2.27%  '

# Row 15: D15, E15
$ws.Range("D15").Value = '21.84'
$ws.Range("E15").Value = '  +5.47%  '

# Row 16: D16
$ws.Range("D16").Value = '3.954.57'

# Row 17: D17, E17
$ws.Range("D17").Value = '14.10'
$ws.Range("E17").Value = '  -1.46%  '

# Row 18: D18, E18
$ws.Range("D18").Value = '1.24'
$ws.Range("E18").Value = '  -2.03%  '

# Row 19: D19, E19
$ws.Range("D19").Value = '72.541.39'
$ws.Range("E19").Value = '  -0.64%  '

# Row 20: E20
$ws.Range("E20").Value = '  -1.01%  '

# Row 21: D21, E21
$ws.Range("D21").Value = '443.94'
$ws.Range("E21").Value = '  +0.02%  '

# Row 22: D22, E22
$ws.Range("D22").Value = '4.76'
$ws.Range("E22").Value = '  -0.03%  '

# Row 23: D23, E23
$ws.Range("D23").Value = '95.70'
$ws.Range("E23").Value = '  -1.91%  '

# Row 24: D24, E24
$ws.Range("D24").Value = '3.32'
$ws.Range("E24").Value = '  -5.85%  '

# Row 25: D25, E25
$ws.Range("D25").Value = '14.15'
$ws.Range("E25").Value = '  -1.44%  '

# Row 26: D26, E26
$ws.Range("D26").Value = '4.27'
$ws.Range("E26").Value = '  -0.77%  '

# Row 27: D27, E27
$ws.Range("D27").Value = '11.13'
$ws.Range("E27").Value = '  -2.87%  '

# Row 28: B28, C28, D28, E28
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").Value = '5.91'
$ws.Range("E28").Value = '  -0.52%  '

# Row 29: B29, C29, D29, E29
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '10.41'
$ws.Range("E29").Value = '  -3.93%  '

# Row 30: D30, E30
$ws.Range("D30").Value = '35.87'
$ws.Range("E30").Value = '  -3.01%  '

# Row 31: D31, E31
$ws.Range("D31").Value = '7.92'
$ws.Range("E31").Value = '  -0.26%  '

# Row 32: D32, E32
$ws.Range("D32").Value = '13.91'
$ws.Range("E32").Value = '  +1.62%  '

# Row 33: D33, E33
$ws.Range("D33").Value = '50.06'
$ws.Range("E33").Value = '  +2.33%  '

# Row 34: E34
$ws.Range("E34").Value = '  -4.00%  '

# Row 35: D35, E35
$ws.Range("D35").Value = '0.0₃0994'
$ws.Range("E35").Value = '  +13.38%  '

# Row 36: D36, E36
$ws.Range("D36").Value = '68.52'
$ws.Range("E36").Value = '  -5.01%  '

# Row 37: D37, E37
$ws.Range("D37").Value = '633.71'
$ws.Range("E37").Value = '  -7.93%  '

# Row 38: D38, E38
$ws.Range("D38").Value = '0.428'
$ws.Range("E38").Value = '  -4.90%  '

# Row 39: D39, E39
$ws.Range("D39").Value = '3.44'
$ws.Range("E39").Value = '  +2.00%  '

# Row 40: E40
$ws.Range("E40").Value = '  +0.19%  '

# Row 41: E41
$ws.Range("E41").Value = '  -2.01%  '

# Row 42: E42
$ws.Range("E42").Value = '  +0.06%  '

# Row 43: D43, E43
$ws.Range("D43").Value = '3.28'
$ws.Range("E43").Value = '  +45.68%  '

# Row 44: B44, C44, D44, E44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0478'
$ws.Range("E44").Value = '  -3.60%  '

# Row 45: B45, C45, D45, E45
$ws.Range("B45").Value = 'THORChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D45").Value = '10.56'
$ws.Range("E45").Value = '  -6.57%  '

# Row 46: E46
$ws.Range("E46").Value = '  -2.22%  '

# Row 47: D47, E47
$ws.Range("D47").Value = '2.61'
$ws.Range("E47").Value = '  -3.85%  '

# Row 48: B48, C48, D48, E48
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '2.87'
$ws.Range("E48").Value = '  -16.64%  '

# Row 49: B49, C49, D49, E49
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '3.37'
$ws.Range("E49").Value = '  -0.89%  '

# Row 50: B50, C50, D50, E50
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '0.000282'
$ws.Range("E50").Value = '  +3.75%  '

# Row 51: D51, E51
$ws.Range("D51").Value = '2.842.00'
$ws.Range("E51").Value = '  +1.37%  '

# Restore column D to General number format (matches source) now that
# the text values are committed, so only the displayed format reverts.
$ws.Range("D2:D51").NumberFormat = "General"
